{"js": "// Convert the two custom \"m:\" fields (fldChar/instrText field codes) into\n// literal M2Doc template text: \" m:'Some value'.setDocumentContentType() \"\n// becomes \"{m:'Some value'.setDocumentContentType()}\" and\n// \" m:''.getDocumentContentType() \" becomes \"{m:''.getDocumentContentType()}\".\n// TokenIteratorFieldRewriterSplit: the parser now reads plain text tokens\n// instead of relying on Word field codes, so the fields are unlinked/removed\n// and replaced by their equivalent literal text, wrapped in \"{ }\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Map of paragraph index -> literal replacement text (braces replace the\n// leading/trailing spaces that used to delimit the field instructions).\nconst replacements = {\n  1: \"{m:'Some value'.setDocumentContentType()}\",\n  2: \"{m:''.getDocumentContentType()}\"\n};\n\nfor (const idxStr of Object.keys(replacements)) {\n  const idx = Number(idxStr);\n  const para = paragraphs.items[idx];\n\n  // Load this paragraph's fields so we can drop the field characters\n  // (fldChar begin/end + instrText runs) before inserting the plain text.\n  const fields = para.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  for (const field of fields.items) {\n    field.delete();\n  }\n  await context.sync();\n\n  // Replace the (now field-free) paragraph content with the literal text.\n  const range = para.getRange(\"Whole\");\n  range.insertText(replacements[idx], \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Convert the two custom \"m:\" fields (fldChar/instrText field codes) into\n# literal M2Doc template text: \" m:'Some value'.setDocumentContentType() \"\n# becomes \"{m:'Some value'.setDocumentContentType()}\" and\n# \" m:''.getDocumentContentType() \" becomes \"{m:''.getDocumentContentType()}\".\n# TokenIteratorFieldRewriterSplit: the parser now reads plain text tokens\n# instead of relying on Word field codes, so the fields are deleted and\n# replaced by their equivalent literal text, wrapped in \"{ }\".\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    2 = \"{m:'Some value'.setDocumentContentType()}\"\n    3 = \"{m:''.getDocumentContentType()}\"\n}\n\nforeach ($idx in $replacements.Keys) {\n    $rng = $d.Paragraphs($idx).Range\n\n    # Drop the field characters (fldChar begin/end + instrText runs) that\n    # make up the field in this paragraph before inserting the plain text.\n    foreach ($f in $rng.Fields) {\n        $f.Delete()\n    }\n\n    $rng.Text = $replacements[$idx]\n}\n"}
